# CW2-member-contribution.xlsx edit
# Commit message: "adds sort by name, impove sorting code"
#
# The underlying change re-scores several tasks to 0 marks awarded (column C)
# and rebalances the "Sort search result" task contribution split between the
# first two students (columns D/E), reflecting sorting-code work being
# re-attributed. All dependent totals/formulas recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Group Size 3")
$ws.Activate()

# Marks Awarded (column C) changes -- several tasks re-scored to 0
$ws.Range("C8").Value  = 0   # Front end HTML & CSS
$ws.Range("C9").Value  = 0   # Product search
$ws.Range("C12").Value = 0   # Registraton and log in/out
$ws.Range("C13").Value = 0   # Editing of customer details.
$ws.Range("C16").Value = 0   # MongoDB dump
$ws.Range("C17").Value = 0   # Report

# Sort search result (row 10) contribution split rebalanced 0/1 -> 0.5/0.5
$ws.Range("D10").Value = 0.5
$ws.Range("E10").Value = 0.5

# Update the active selection to match the edited workbook (C18, no frozen
# top-left scroll position)
$ws.Range("C18").Select()
